# Commit: "Added new temp class for creating 200 jobs. Added new functions and TCs."
# This populates the previously-empty rows 59-64 on the "TestData" sheet with new
# Schedule-Job related key/value test data (ScheduleJobName, ScheduleJobDescription,
# ScheduleJobDate, ScheduleJobTime, ScheduleJobTimeZone, ScheduleJobRecurrence),
# continuing directly after the existing "TC018 - Send Dashboard via email" block
# (row 58 = AppendTimeStampHelpText).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Carry over the look & feel (font/alignment/row-style) of the last populated
# detail row (58) onto the new rows before writing the new values into them.
$ws.Range("A58:B58").Copy()
$ws.Range("A59:B64").PasteSpecial(-4122)

$values = @(
    @("ScheduleJobName",        "UpdateJob"),
    @("ScheduleJobDescription", "Automation_Job_Name"),
    @("ScheduleJobDate",        "10/31/2021"),
    @("ScheduleJobTime",        "03:30 PM"),
    @("ScheduleJobTimeZone",    "GMT+03:00"),
    @("ScheduleJobRecurrence",  "No Recurrence")
)

$row = 59
foreach ($pair in $values) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}
